# Daily TGP (terminal gate pricing) roll-forward:
# - New day's rows (serial 46051) written on top of each 2-3 row date group
# - Previous top rows (old 46050) shift down to where 46046 used to be
# - Rows that used to hold 46046 are overwritten by the shifted-down 46050 data
# (46046 data itself is retired / drops off the sheet)
#
# Columns: A = Effective Date (serial), D = Diesel, E = ULP, F = PULP, G = e10
# Blank / unused cells in a group (merged "N/A" text cells, or groups with no
# G column) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row=8;  A=46051; D=158.25; E=149.89; F=159.89; G=149.91 },
    @{ Row=9;  A=46051; D=158.25; E=149.89; F=159.89; G=149.91 },
    @{ Row=10; A=46051; D=159.12; E=151.73; F=161.73; G=152.13 },
    @{ Row=11; A=46050; D=158.59; E=150.49; F=160.49; G=150.51 },
    @{ Row=12; A=46050; D=158.59; E=150.49; F=160.49; G=150.51 },
    @{ Row=13; A=46050; D=159.45; E=152.33; F=162.33; G=152.73 },
    @{ Row=17; A=46051; D=162.61; E=154.32; F=164.32 },
    @{ Row=18; A=46050; D=162.91; E=154.89; F=164.89 },
    @{ Row=22; A=46051; D=159.21; E=151.46; F=161.06; G=152.54 },
    @{ Row=23; A=46051; D=163.89; E=157.07; F=167.07 },
    @{ Row=24; A=46051; D=164.04; E=157.75; F=167.75 },
    @{ Row=25; A=46051; D=164.03; E=157.29; F=167.29; G=157.41 },
    @{ Row=26; A=46051; D=163.62; E=158.89; F=168.89 },
    @{ Row=27; A=46050; D=159.67; E=152.06; F=161.66; G=153.14 },
    @{ Row=28; A=46050; D=164.23; E=157.66; F=167.66 },
    @{ Row=29; A=46050; D=164.37; E=158.36; F=168.36 },
    @{ Row=30; A=46050; D=164.36; E=157.9;  F=167.9;  G=158.02 },
    @{ Row=31; A=46050; D=163.95; E=159.51; F=169.51 },
    @{ Row=35; A=46051; D=158.02; E=148.74; F=157.74 },
    @{ Row=36; A=46050; D=158.35; E=149.33; F=158.33 },
    @{ Row=40; A=46051; D=163.49; E=156.32; F=166.32 },
    @{ Row=41; A=46051; D=163.21; E=156.74; F=166.74 },
    @{ Row=42; A=46050; D=163.85; E=157.09; F=167.09 },
    @{ Row=43; A=46050; D=163.57; E=157.51; F=167.51 },
    @{ Row=47; A=46051; D=157.86; E=150.44; F=160.44 },
    @{ Row=48; A=46051; D=157.45; E=150.36; F=160.36 },
    @{ Row=49; A=46050; D=157.47; E=150.83; F=160.83 },
    @{ Row=50; A=46050; D=157.06; E=150.74; F=160.74 },
    @{ Row=54; A=46051; D=172.51; E=164.78; F=174.78 },
    @{ Row=55; A=46051; D=165.19; E=163.01; F=173.01 },
    @{ Row=56; A=46051; D=162.06 },
    @{ Row=57; A=46051; D=162.52; E=157.44 },
    @{ Row=58; A=46051; D=158.29; E=153.33; F=163.33 },
    @{ Row=59; A=46051; D=164.95; E=162.89 },
    @{ Row=60; A=46050; D=172.87; E=165.47; F=175.47 },
    @{ Row=61; A=46050; D=165.58; E=163.58; F=173.58 },
    @{ Row=62; A=46050; D=162.51 },
    @{ Row=63; A=46050; D=162.93; E=158 },
    @{ Row=64; A=46050; D=158.7;  E=153.9;  F=163.9 },
    @{ Row=65; A=46050; D=165.24; E=163.55 }
)

foreach ($u in $updates) {
    $r = $u.Row
    $ws.Cells.Item($r, 1).Value = $u.A          # A: Effective Date
    if ($u.ContainsKey('D')) { $ws.Cells.Item($r, 4).Value = $u.D }   # D: Diesel
    if ($u.ContainsKey('E')) { $ws.Cells.Item($r, 5).Value = $u.E }   # E: ULP
    if ($u.ContainsKey('F')) { $ws.Cells.Item($r, 6).Value = $u.F }   # F: PULP
    if ($u.ContainsKey('G')) { $ws.Cells.Item($r, 7).Value = $u.G }   # G: e10
}
